$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 ("I0") and J1 ("IF") ---------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the look of the existing header cells (bold, centered, thin box)
# by copying the format from H1 (same header style) onto the new cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows: add I/J values -----------------------------------------
# Row 2 is special (I0=9, IF=9); every other row repeats its IP (H) value
# into IF (J) and uses 1 for I0.
$newCols = @{
    2  = @(9, 9)
    3  = @(1, 5)
    4  = @(1, 4)
    5  = @(1, 5)
    6  = @(1, 5)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 6)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 4)
    13 = @(1, 4)
    14 = @(1, 4)
    15 = @(1, 5)
    16 = @(1, 3)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 4)
    21 = @(1, 4)
}

foreach ($r in 2..21) {
    $vals = $newCols[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
